$wb = $excel.ActiveWorkbook

# "Logs" is the active sheet and holds the mail log table.
$logs = $wb.ActiveSheet

# Append the new test-mail entry as row 51 (table currently ends at row 50).
$logs.Range("A51").Value = "Laat maar weten of er nieuws is"
$logs.Range("B51").Value = "mailmind.test@zohomail.eu"
$logs.Range("C51").Value = "Testmail #10: Laat maar weten of er nieuws is"
$logs.Range("D51").Value = "Overig"
$logs.Range("E51").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Range("F51").Value = "2025-08-05 19:42:19"
$logs.Range("G51").Value = "Ja"
$logs.Range("H51").Value = "Ja"
$logs.Range("I51").Value = "Nee"
$logs.Range("J51").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered,
# same way Excel does when a table/range grows by one row.
$ranges = @("D2:D50", "G2:G50", "H2:H50", "I2:I50", "J2:J50")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newRange = $logs.Range("$($col)2:$($col)51")
    $fcs = $logs.Range($addr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# "Dashboard" tallies the count of log entries per category; "Overig" grew by one.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 9
